$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Terraform, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Collapse(0)
$rng.InsertAfter("Packer, ")
